{"js": "const replacements = [\n  [\"912\u00f76=152, 0\", \"108\u00f77=15, 3\"],\n  [\"308\u00f75=61, 3\", \"337\u00f78=42, 1\"],\n  [\"713\u00f77=101, 6\", \"774\u00f72=387, 0\"],\n  [\"424\u00f79=47, 1\", \"978\u00f78=122, 2\"],\n  [\"779\u00f77=111, 2\", \"541\u00f74=135, 1\"],\n  [\"785\u00f74=196, 1\", \"453\u00f72=226, 1\"],\n  [\"428\u00f73=142, 2\", \"860\u00f75=172, 0\"],\n  [\"901\u00f78=112, 5\", \"898\u00f72=449, 0\"],\n  [\"218\u00f79=24, 2\", \"750\u00f73=250, 0\"],\n  [\"208\u00f74=52, 0\", \"588\u00f77=84, 0\"],\n  [\"374\u00f72=187, 0\", \"765\u00f74=191, 1\"],\n  [\"109\u00f79=12, 1\", \"388\u00f79=43, 1\"],\n  [\"604\u00f72=302, 0\", \"112\u00f72=56, 0\"],\n  [\"598\u00f73=199, 1\", \"584\u00f78=73, 0\"],\n  [\"998\u00f72=499, 0\", \"957\u00f72=478, 1\"],\n  [\"889\u00f77=127, 0\", \"534\u00f79=59, 3\"],\n  [\"838\u00f72=419, 0\", \"621\u00f78=77, 5\"],\n  [\"908\u00f78=113, 4\", \"671\u00f78=83, 7\"],\n  [\"856\u00f78=107, 0\", \"892\u00f74=223, 0\"],\n  [\"494\u00f76=82, 2\", \"158\u00f76=26, 2\"],\n  [\"266\u00f75=53, 1\", \"325\u00f77=46, 3\"],\n  [\"422\u00f77=60, 2\", \"373\u00f76=62, 1\"],\n  [\"466\u00f75=93, 1\", \"307\u00f77=43, 6\"],\n  [\"168\u00f75=33, 3\", \"401\u00f74=100, 1\"],\n  [\"569\u00f76=94, 5\", \"731\u00f77=104, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"912\u00f76=152, 0\"; New = \"108\u00f77=15, 3\" },\n    @{ Old = \"308\u00f75=61, 3\"; New = \"337\u00f78=42, 1\" },\n    @{ Old = \"713\u00f77=101, 6\"; New = \"774\u00f72=387, 0\" },\n    @{ Old = \"424\u00f79=47, 1\"; New = \"978\u00f78=122, 2\" },\n    @{ Old = \"779\u00f77=111, 2\"; New = \"541\u00f74=135, 1\" },\n    @{ Old = \"785\u00f74=196, 1\"; New = \"453\u00f72=226, 1\" },\n    @{ Old = \"428\u00f73=142, 2\"; New = \"860\u00f75=172, 0\" },\n    @{ Old = \"901\u00f78=112, 5\"; New = \"898\u00f72=449, 0\" },\n    @{ Old = \"218\u00f79=24, 2\"; New = \"750\u00f73=250, 0\" },\n    @{ Old = \"208\u00f74=52, 0\"; New = \"588\u00f77=84, 0\" },\n    @{ Old = \"374\u00f72=187, 0\"; New = \"765\u00f74=191, 1\" },\n    @{ Old = \"109\u00f79=12, 1\"; New = \"388\u00f79=43, 1\" },\n    @{ Old = \"604\u00f72=302, 0\"; New = \"112\u00f72=56, 0\" },\n    @{ Old = \"598\u00f73=199, 1\"; New = \"584\u00f78=73, 0\" },\n    @{ Old = \"998\u00f72=499, 0\"; New = \"957\u00f72=478, 1\" },\n    @{ Old = \"889\u00f77=127, 0\"; New = \"534\u00f79=59, 3\" },\n    @{ Old = \"838\u00f72=419, 0\"; New = \"621\u00f78=77, 5\" },\n    @{ Old = \"908\u00f78=113, 4\"; New = \"671\u00f78=83, 7\" },\n    @{ Old = \"856\u00f78=107, 0\"; New = \"892\u00f74=223, 0\" },\n    @{ Old = \"494\u00f76=82, 2\"; New = \"158\u00f76=26, 2\" },\n    @{ Old = \"266\u00f75=53, 1\"; New = \"325\u00f77=46, 3\" },\n    @{ Old = \"422\u00f77=60, 2\"; New = \"373\u00f76=62, 1\" },\n    @{ Old = \"466\u00f75=93, 1\"; New = \"307\u00f77=43, 6\" },\n    @{ Old = \"168\u00f75=33, 3\"; New = \"401\u00f74=100, 1\" },\n    @{ Old = \"569\u00f76=94, 5\"; New = \"731\u00f77=104, 3\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
